$wb = $excel.ActiveWorkbook

# Rename the "DropdownLists" sheet to "Dropdownlists"
$ws2 = $wb.Worksheets.Item("DropdownLists")
$ws2.Name = "Dropdownlists"

# Make the renamed "Dropdownlists" sheet the active sheet/tab
$ws2.Activate()

$wb.Save()
